# Apply updated betting-odds values to Sheet1, as described by the source diff.
# Each assignment below sets a single cell to its new numeric value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("G6").Value = 2.3
$ws.Range("I6").Value = 3.75
$ws.Range("J6").Value = 3.2
$ws.Range("L6").Value = 4.5
$ws.Range("R6").Value = 1.33
$ws.Range("W6").Value = 5.5
$ws.Range("AA6").Value = 26
$ws.Range("AO6").Value = 15
$ws.Range("BB6").Value = 451

# Row 7
$ws.Range("G7").Value = 1.9
$ws.Range("H7").Value = 2.9
$ws.Range("I7").Value = 5.5
$ws.Range("L7").Value = 6
$ws.Range("W7").Value = 4.5
$ws.Range("AJ7").Value = 21
$ws.Range("AK7").Value = 67
$ws.Range("AO7").Value = 11
$ws.Range("AZ7").Value = 151

# Row 8
$ws.Range("G8").Value = 2.55
$ws.Range("H8").Value = 2.7
$ws.Range("I8").Value = 3.4
$ws.Range("J8").Value = 3.5
$ws.Range("K8").Value = 1.8
$ws.Range("M8").Value = 1.17
$ws.Range("N8").Value = 5
$ws.Range("O8").Value = 1.67
$ws.Range("P8").Value = 2.1
$ws.Range("Q8").Value = 3.4
$ws.Range("R8").Value = 1.33
$ws.Range("S8").Value = 1.73
$ws.Range("T8").Value = 2
$ws.Range("X8").Value = 10
$ws.Range("Y8").Value = 12
$ws.Range("Z8").Value = 26
$ws.Range("AA8").Value = 29
$ws.Range("AB8").Value = 51
$ws.Range("AF8").Value = 101
$ws.Range("AH8").Value = 6.5
$ws.Range("AN8").Value = 4.33
$ws.Range("AO8").Value = 17
$ws.Range("AT8").Value = 2
$ws.Range("AV8").Value = 101
$ws.Range("AX8").Value = 21
$ws.Range("BB8").Value = 501

# Row 22
$ws.Range("N22").Value = 8
$ws.Range("Q22").Value = 2.3
$ws.Range("R22").Value = 1.6

# Row 23
$ws.Range("G23").Value = 2.15
$ws.Range("H23").Value = 3.1
$ws.Range("I23").Value = 3.8
$ws.Range("J23").Value = 3
$ws.Range("L23").Value = 4.5
$ws.Range("S23").Value = 1.57
$ws.Range("T23").Value = 2.25
$ws.Range("W23").Value = 6
$ws.Range("X23").Value = 9
$ws.Range("Z23").Value = 19
$ws.Range("AA23").Value = 21
$ws.Range("AH23").Value = 8.5
$ws.Range("AI23").Value = 17
$ws.Range("AJ23").Value = 13
$ws.Range("AL23").Value = 34
$ws.Range("AM23").Value = 41
$ws.Range("AN23").Value = 4
$ws.Range("AO23").Value = 13
$ws.Range("AP23").Value = 29
$ws.Range("AR23").Value = 81
$ws.Range("AT23").Value = 2.25
$ws.Range("AY23").Value = 34

# Row 24
$ws.Range("M24").Value = 1.03
$ws.Range("N24").Value = 15
$ws.Range("Q24").Value = 1.67
$ws.Range("R24").Value = 2.15
$ws.Range("S24").Value = 1.3
$ws.Range("T24").Value = 3.4
$ws.Range("U24").Value = 2
$ws.Range("V24").Value = 1.75
$ws.Range("AB24").Value = 26
$ws.Range("AC24").Value = 15
$ws.Range("AG24").Value = 401
$ws.Range("AP24").Value = 17
$ws.Range("AQ24").Value = 17
$ws.Range("AT24").Value = 3.4
$ws.Range("AU24").Value = 9

# Row 25
$ws.Range("G25").Value = 1.91
$ws.Range("H25").Value = 3.6
$ws.Range("I25").Value = 4
$ws.Range("J25").Value = 2.6
$ws.Range("Q25").Value = 2.07
$ws.Range("R25").Value = 1.83
$ws.Range("S25").Value = 1.4
$ws.Range("T25").Value = 2.75
$ws.Range("Y25").Value = 8.5
$ws.Range("AA25").Value = 15
$ws.Range("AC25").Value = 10
$ws.Range("AE25").Value = 15
$ws.Range("AH25").Value = 11
$ws.Range("AO25").Value = 10
$ws.Range("AQ25").Value = 34
$ws.Range("AT25").Value = 2.75
$ws.Range("AW25").Value = 6
$ws.Range("BB25").Value = 201

# Row 26
$ws.Range("G26").Value = 2.2
$ws.Range("H26").Value = 2.9
$ws.Range("J26").Value = 3.2
$ws.Range("O26").Value = 1.67
$ws.Range("P26").Value = 2.1
$ws.Range("Q26").Value = 3.4
$ws.Range("R26").Value = 1.33
$ws.Range("S26").Value = 1.73
$ws.Range("T26").Value = 2.08
$ws.Range("Y26").Value = 9

# Row 27
$ws.Range("G27").Value = 2
$ws.Range("L27").Value = 5
$ws.Range("M27").Value = 1.13
$ws.Range("N27").Value = 6
$ws.Range("O27").Value = 1.53
$ws.Range("P27").Value = 2.38
$ws.Range("Q27").Value = 2.7
$ws.Range("R27").Value = 1.44
$ws.Range("AE27").Value = 19
$ws.Range("AH27").Value = 9
$ws.Range("AI27").Value = 21
$ws.Range("AU27").Value = 9.5

# Row 34
$ws.Range("G34").Value = 1.85
$ws.Range("H34").Value = 3.25
$ws.Range("I34").Value = 4.5
$ws.Range("J34").Value = 2.63
$ws.Range("K34").Value = 1.95
$ws.Range("L34").Value = 5
$ws.Range("M34").Value = 1.11
$ws.Range("N34").Value = 6.5
$ws.Range("Q34").Value = 2.5
$ws.Range("R34").Value = 1.5
$ws.Range("U34").Value = 2.2
$ws.Range("V34").Value = 1.62
$ws.Range("W34").Value = 5.5
$ws.Range("X34").Value = 7.5
$ws.Range("Y34").Value = 9.5
$ws.Range("Z34").Value = 15
$ws.Range("AA34").Value = 19
$ws.Range("AC34").Value = 6.5
$ws.Range("AD34").Value = 6.5
$ws.Range("AH34").Value = 9.5
$ws.Range("AI34").Value = 21
$ws.Range("AJ34").Value = 17
$ws.Range("AL34").Value = 41
$ws.Range("AN34").Value = 3.6
$ws.Range("AO34").Value = 11
$ws.Range("AQ34").Value = 41
$ws.Range("AW34").Value = 6
$ws.Range("AX34").Value = 29
$ws.Range("AZ34").Value = 101

# Row 42
$ws.Range("G42").Value = 1.8
$ws.Range("I42").Value = 5.5
$ws.Range("L42").Value = 5.5
$ws.Range("O42").Value = 1.5
$ws.Range("P42").Value = 2.5
$ws.Range("AK42").Value = 51
$ws.Range("AM42").Value = 51
$ws.Range("AN42").Value = 3.6
$ws.Range("AX42").Value = 29
$ws.Range("BA42").Value = 151

# Row 62
$ws.Range("U62").Value = 1.87
$ws.Range("V62").Value = 1.87

# Row 63
$ws.Range("G63").Value = 4.75
$ws.Range("I63").Value = 1.7
$ws.Range("J63").Value = 4.5
$ws.Range("K63").Value = 2.4
$ws.Range("L63").Value = 2.25
$ws.Range("Q63").Value = 1.62
$ws.Range("R63").Value = 2.25
$ws.Range("U63").Value = 1.62
$ws.Range("V63").Value = 2.2
$ws.Range("AB63").Value = 34
$ws.Range("AC63").Value = 15
$ws.Range("AH63").Value = 9
$ws.Range("AI63").Value = 9.5
$ws.Range("AL63").Value = 13
$ws.Range("AP63").Value = 26
$ws.Range("AQ63").Value = 67
$ws.Range("AV63").Value = 41
$ws.Range("AW63").Value = 4
$ws.Range("BC63").Value = 451

# Row 64
$ws.Range("H64").Value = 3.9
$ws.Range("I64").Value = 5
$ws.Range("J64").Value = 2.12
$ws.Range("K64").Value = 2.25
$ws.Range("Q64").Value = 1.65
$ws.Range("R64").Value = 2
$ws.Range("X64").Value = 8
$ws.Range("Z64").Value = 12
$ws.Range("AC64").Value = 12.5
$ws.Range("AD64").Value = 7.8
$ws.Range("AE64").Value = 15
$ws.Range("AF64").Value = 60
$ws.Range("AJ64").Value = 15.5
$ws.Range("AN64").Value = 3.5
$ws.Range("AX64").Value = 27
$ws.Range("AY64").Value = 30

# Row 80
$ws.Range("G80").Value = 1.62
$ws.Range("O80").Value = 1.4
$ws.Range("P80").Value = 2.75
$ws.Range("U80").Value = 2.2
$ws.Range("V80").Value = 1.62
$ws.Range("AA80").Value = 15
$ws.Range("AB80").Value = 34
$ws.Range("AS80").Value = 201
$ws.Range("AU80").Value = 9.5
$ws.Range("AZ80").Value = 126
$ws.Range("BA80").Value = 151

# Row 81
$ws.Range("Q81").Value = 1.62
$ws.Range("R81").Value = 2.25

